$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the column headers: "<name>_old" -> "<name>_FV2404",
#    "<name>_new" -> "<name>_FV2410" (the "diff" header stays unchanged).
# ---------------------------------------------------------------------------
$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------------
# 2) Turn A1:U91 into a real Excel Table (ListObject) named "Table1".
#    The header row (row 1) already carries explicit direct formatting
#    (bold/fill/border) from the original template. Building the table
#    directly on top of it would make Excel "bake" that pre-existing
#    formatting into a headerRowDxfId/dxf - something the target workbook
#    does not have. So the table is first created on a blank, unformatted
#    staging area (with the very same header captions) and then resized
#    onto the real A1:U91 range, which keeps the existing cell styles
#    untouched and avoids generating any extra dxf.
# ---------------------------------------------------------------------------
$stageRow = $ws.UsedRange.Rows.Count + 10

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item($stageRow, $i + 1).Value = $headers[$i]
    $ws.Cells.Item($stageRow + 1, $i + 1).Value = "x"
}

$stageRange = $ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow + 1, 21))
$lo = $ws.ListObjects.Add(1, $stageRange, $null, 1)
$lo.TableStyle = $null
$lo.Resize($ws.Range("A1:U91"))

# Clean up the temporary staging rows again.
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow + 1, 21)).Clear()

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split beneath row 1, bottom-left pane active).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
